$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the
# three new header cells so they pick up the same bold/border/centered
# style used by the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Set the new header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 86  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 76  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
